$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("DBD")

# ---------------------------------------------------------------------------
# Two new rows (PfxPath / PfxAuth) are being inserted into the field list at
# rows 78-79. The SEQ column (A) already holds a simple running sequence
# (row number - 8) for every populated row, so it does not need to be
# shifted - only columns B:H (field definitions) move down by two rows for
# the previously-existing CreateDate/CreateEmpNo/LastUpdate/LastUpdateEmpNo
# rows. Column H is handled separately below because only rows 77-79 use it.
# ---------------------------------------------------------------------------

# 1) Shift B:G content down by 2 rows: 78->80, 79->81, 80->82, 81->83
#    (copy starting from the bottom so sources are not clobbered)
$ws.Range("B81:G81").Copy($ws.Range("B83:G83"))
$ws.Range("B80:G80").Copy($ws.Range("B82:G82"))
$ws.Range("B79:G79").Copy($ws.Range("B81:G81"))
$ws.Range("B78:G78").Copy($ws.Range("B80:G80"))
$excel.CutCopyMode = $false

# The custom row height (22.2) that belonged to the old "CreateEmpNo" row
# (79) needs to move along with that row's content, to its new home at 81.
$ws.Rows(79).AutoFit()
$ws.Rows(81).RowHeight = 22.2

# Fill in the SEQ numbers for the two rows that now contain data that used
# to be blank (the former empty rows 82/83 receive LastUpdate* content).
$ws.Range("A82").Value = 74
$ws.Range("A83").Value = 75

# 2) Apply the "recently added" (red font) styling used for new entries to
#    the two brand-new rows, copied from the previous new-entry row (77).
$ws.Range("B77:G77").Copy()
$ws.Range("B78:G79").PasteSpecial(-4122)
$ws.Range("H77").Copy()
$ws.Range("H78:H79").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Add the remark to row 77 (L6972Flag), which was previously blank.
$ws.Range("H77").Value = "2023/6/13智偉新增"

# 4) Populate the two new field-definition rows.
$ws.Range("B78").Value = "PfxPath"
$ws.Range("C78").Value = "憑證路徑"
$ws.Range("D78").Value = "VARCHAR2"
$ws.Range("E78").Value = 100
$ws.Range("H78").Value = "2023/6/19智偉新增"

$ws.Range("B79").Value = "PfxAuth"
$ws.Range("C79").Value = "憑證認證"
$ws.Range("D79").Value = "VARCHAR2"
$ws.Range("E79").Value = 100
$ws.Range("H79").Value = "2023/6/19智偉新增"

# 5) Append two new blank rows at the bottom (85-86), matching the style of
#    the existing trailing blank row (84).
$ws.Range("A84:G84").Copy($ws.Range("A85:G85"))
$ws.Range("A84:G84").Copy($ws.Range("A86:G86"))
$excel.CutCopyMode = $false

# 6) Update the view to reflect the current working position (scrolled down
#    a couple more rows, selection moved to the new blank row being edited).
$ws.Range("A66").Select()
$excel.ActiveWindow.ScrollRow = 66
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C82").Select()
